# Refresh the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions job, see commit message). Rows 32-34 also got
# reordered (Toncoin now ranks above InjectiveProtocol/RenderToken).
#
# Price/volume cells are plain text in this sheet (e.g. "52.337.30",
# "  -0.08%  "), so a helper writes the value and, for anything that
# Excel's general-number parser would otherwise coerce into a float
# (dropping formatting / introducing rounding noise), forces it back to
# literal text the same way typing a leading apostrophe would in the UI.
function Set-CellText($sheet, $ref, $text) {
    $sheet.Range($ref).Value = $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "52.337.30"
Set-CellText $ws "E2" "  -0.08%  "
Set-CellText $ws "D3" "2.934.33"
Set-CellText $ws "E3" "  +0.71%  "
Set-CellText $ws "E4" "  -0.04%  "
Set-CellText $ws "D5" "'357.93"
Set-CellText $ws "E5" "  +1.38%  "
Set-CellText $ws "D6" "'109.97"
Set-CellText $ws "E6" "  -2.42%  "
Set-CellText $ws "D7" "'0.569"
Set-CellText $ws "E7" "  +1.34%  "
Set-CellText $ws "E8" "  +0.03%  "
Set-CellText $ws "D9" "'0.633"
Set-CellText $ws "E9" "  +0.23%  "
Set-CellText $ws "D10" "'39.13"
Set-CellText $ws "E10" "  -2.37%  "
Set-CellText $ws "E11" "  +1.46%  "
Set-CellText $ws "D12" "'0.0873"
Set-CellText $ws "E12" "  +0.54%  "
Set-CellText $ws "D13" "'19.58"
Set-CellText $ws "E13" "  -1.67%  "
Set-CellText $ws "D14" "'7.78"
Set-CellText $ws "E14" "  -0.54%  "
Set-CellText $ws "D15" "3.397.61"
Set-CellText $ws "E15" "  +0.72%  "
Set-CellText $ws "D16" "2.947.44"
Set-CellText $ws "E16" "  +0.37%  "
Set-CellText $ws "D17" "'0.988"
Set-CellText $ws "E17" "  -1.80%  "
Set-CellText $ws "D18" "52.328.65"
Set-CellText $ws "E18" "  -0.14%  "
Set-CellText $ws "D19" "'3.52"
Set-CellText $ws "E19" "  +5.87%  "
Set-CellText $ws "D20" "'7.59"
Set-CellText $ws "E20" "  -0.77%  "
Set-CellText $ws "D21" "'13.97"
Set-CellText $ws "E21" "  -1.81%  "
Set-CellText $ws "D22" "0.0₃0985"
Set-CellText $ws "E22" "  +0.24%  "
Set-CellText $ws "D23" "'70.57"
Set-CellText $ws "E23" "  -0.59%  "
Set-CellText $ws "D24" "'269.01"
Set-CellText $ws "E24" "  -0.49%  "
Set-CellText $ws "D25" "'2.81"
Set-CellText $ws "E25" "  +0.77%  "
Set-CellText $ws "D26" "'0.185"
Set-CellText $ws "E26" "  +7.06%  "
Set-CellText $ws "D27" "'7.79"
Set-CellText $ws "E27" "  +16.66%  "
Set-CellText $ws "D28" "'27.09"
Set-CellText $ws "E28" "  +0.97%  "
Set-CellText $ws "E29" "  +0.11%  "
Set-CellText $ws "E30" "  +7.85%  "
Set-CellText $ws "D31" "'10.52"
Set-CellText $ws "B32" "Toncoin"
Set-CellText $ws "C32" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D32" "'2.29"
Set-CellText $ws "E32" "  +1.41%  "
Set-CellText $ws "B33" "InjectiveProtocol"
Set-CellText $ws "C33" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText $ws "D33" "'37.61"
Set-CellText $ws "E33" "  -1.34%  "
Set-CellText $ws "B34" "RenderToken"
Set-CellText $ws "C34" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D34" "'6.20"
Set-CellText $ws "E34" "  -2.74%  "
Set-CellText $ws "D35" "'52.33"
Set-CellText $ws "E35" "  -1.69%  "
Set-CellText $ws "D36" "'0.0445"
Set-CellText $ws "E36" "  -1.70%  "
Set-CellText $ws "E37" "  +0.00%  "
Set-CellText $ws "E38" "  -4.26%  "
Set-CellText $ws "E39" "  -3.19%  "
Set-CellText $ws "D40" "'2.78"
Set-CellText $ws "E41" "  -3.58%  "
Set-CellText $ws "E42" "  +2.41%  "
Set-CellText $ws "D43" "'23.26"
Set-CellText $ws "E43" "  -0.41%  "
Set-CellText $ws "D44" "'119.86"
Set-CellText $ws "E44" "  -1.17%  "
Set-CellText $ws "D45" "'2.17"
Set-CellText $ws "E45" "  -1.06%  "
Set-CellText $ws "E46" "  -2.77%  "
Set-CellText $ws "D47" "'2.48"
Set-CellText $ws "E47" "  -5.19%  "
Set-CellText $ws "D48" "2.132.90"
Set-CellText $ws "E48" "  -3.08%  "
Set-CellText $ws "D49" "'0.252"
Set-CellText $ws "E49" "  -5.17%  "
Set-CellText $ws "D50" "'0.0351"
Set-CellText $ws "E50" "  +2.05%  "
Set-CellText $ws "D51" "'0.932"
Set-CellText $ws "E51" "  -3.81%  "
